$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tests")
$ws2 = $wb.Worksheets.Item("Result")

# ---------------------------------------------------------------------------
# Step 1: Write the brand new text values in the precise order they must
# first appear so that the workbook's shared-string table assigns the same
# indices (17-23) that the target file uses.
# ---------------------------------------------------------------------------
$ws1.Range("A6").Value = "Test_Framework\Test_AppEmail.xaml"                          # new shared string 17
$ws1.Range("E6").Value = "Check Dev exchange folders"                                  # new shared string 18
$ws1.Range("E7").Value = "Check testing exchange folders"                              # new shared string 19
$ws1.Range("E8").Value = "Check production account and exchange folders"               # new shared string 20
$ws2.Range("C6").Value = "FAIL"                                                        # new shared string 21
$ws1.Range("A7").Value = "\Test_Framework\Test_AppEmail.xaml"                          # new shared string 22
$ws2.Range("E6").Value = "AppEx: Folder not found. at Source: Get Exchange Robot Invalid"  # new shared string 23

# ---------------------------------------------------------------------------
# Step 2: Fill in the rest of the new rows (these all reuse existing shared
# strings, so order no longer matters for string-table indices).
# ---------------------------------------------------------------------------

# --- Sheet "Tests" row 6 ---
$ws1.Range("B6").Value = "Success"
$ws1.Range("C6").WrapText = $true
$ws1.Range("D6").Value = "Data\ConfigDev.xlsx"
$ws1.Range("D6").WrapText = $true

# --- Sheet "Tests" row 7 ---
$ws1.Range("B7").Value = "Success"
$ws1.Range("D7").Value = "Data\ConfigUat.xlsx"
$ws1.Range("D7").WrapText = $true

# --- Sheet "Tests" row 8 ---
$ws1.Range("A8").Value = "\Test_Framework\Test_AppEmail.xaml"
$ws1.Range("B8").Value = "Success"
$ws1.Range("D8").Value = "Data\ConfigPrd.xlsx"
$ws1.Range("D8").WrapText = $true

# --- Sheet "Result" row 6 ---
# Columns A:E on this sheet default to the word-wrap style, so cells that
# must stay on the plain "Normal" style (A, B, E) need it applied explicitly.
$ws2.Range("A6").Value = "Test_Framework\Test_AppEmail.xaml"
$ws2.Range("A6").Style = "Normal"
$ws2.Range("B6").Value = "Success"
$ws2.Range("B6").Style = "Normal"
$ws2.Range("C6").WrapText = $true
$ws2.Range("D6").Value = "Data\ConfigDev.xlsx"
$ws2.Range("D6").WrapText = $true
$ws2.Range("E6").Style = "Normal"

# --- Sheet "Result" row 7 ---
$ws2.Range("A7").Value = "\Test_Framework\Test_AppEmail.xaml"
$ws2.Range("A7").Style = "Normal"
$ws2.Range("B7").Value = "Success"
$ws2.Range("B7").Style = "Normal"
$ws2.Range("C7").Style = "Normal"
$ws2.Range("D7").Value = "Data\ConfigUat.xlsx"
$ws2.Range("D7").WrapText = $true
$ws2.Range("E7").Value = "Check testing exchange folders"
$ws2.Range("E7").Style = "Normal"

# --- Sheet "Result" row 8 ---
$ws2.Range("A8").Value = "\Test_Framework\Test_AppEmail.xaml"
$ws2.Range("A8").Style = "Normal"
$ws2.Range("B8").Value = "Success"
$ws2.Range("B8").Style = "Normal"
$ws2.Range("C8").Style = "Normal"
$ws2.Range("D8").Value = "Data\ConfigPrd.xlsx"
$ws2.Range("D8").WrapText = $true
$ws2.Range("E8").Value = "Check production account and exchange folders"
$ws2.Range("E8").Style = "Normal"

# --- Sheet "Result" row 9 : blank row touched but left empty ---
$ws2.Range("A9").Style = "Normal"
$ws2.Range("B9").Style = "Normal"
$ws2.Range("D9").Style = "Normal"
$ws2.Range("E9").Style = "Normal"

# ---------------------------------------------------------------------------
# Step 3: Column D on "Result" grew wider to fit the new, longer values.
# ---------------------------------------------------------------------------
$ws2.Columns.Item(4).ColumnWidth = 31.5

# ---------------------------------------------------------------------------
# Step 4: Restore the view/selection state recorded in the saved workbook.
# ---------------------------------------------------------------------------
$ws1.Activate() | Out-Null
$ws1.Range("A9").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("A9").Select() | Out-Null
